# This workbook holds a weekly/daily price table for Pera (pear) at
# "Vega Monumental Concepción". The edit prepends two new records
# (Packham's Triumph - Primera / Segunda, dated 2023-10-19) at row 740,
# pushing the existing rows 740:859 down to 742:861 (dimension grows
# from A1:T859 to A1:T861).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 740:741; Excel shifts rows 740:859 down to
# 742:861 and copies the row-740 formatting (incl. the date style on
# column D) into the freshly inserted rows.
$ws.Rows("740:741").Insert()

# --- New row 740: Packham's Triumph / Primera -----------------------
$ws.Cells.Item(740, 1).Value = 11
$ws.Cells.Item(740, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(740, 3).Value = "Bíobío"
$ws.Cells.Item(740, 4).Value = 45218
$ws.Cells.Item(740, 5).Value = 8
$ws.Cells.Item(740, 6).Value = "Fruta"
$ws.Cells.Item(740, 7).Value = 100104
$ws.Cells.Item(740, 8).Value = "Frutos de pepita"
$ws.Cells.Item(740, 9).Value = 100104005
$ws.Cells.Item(740, 10).Value = "Pera"
$ws.Cells.Item(740, 11).Value = "Packham's Triumph"
$ws.Cells.Item(740, 12).Value = "Primera"
$ws.Cells.Item(740, 13).Value = 150
$ws.Cells.Item(740, 14).Value = 13000
$ws.Cells.Item(740, 15).Value = 13000
$ws.Cells.Item(740, 16).Value = 13000
$ws.Cells.Item(740, 17).Value = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(740, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(740, 19).Value = 812
$ws.Cells.Item(740, 20).Value = 16

# --- New row 741: Packham's Triumph / Segunda ------------------------
$ws.Cells.Item(741, 1).Value = 11
$ws.Cells.Item(741, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(741, 3).Value = "Bíobío"
$ws.Cells.Item(741, 4).Value = 45218
$ws.Cells.Item(741, 5).Value = 8
$ws.Cells.Item(741, 6).Value = "Fruta"
$ws.Cells.Item(741, 7).Value = 100104
$ws.Cells.Item(741, 8).Value = "Frutos de pepita"
$ws.Cells.Item(741, 9).Value = 100104005
$ws.Cells.Item(741, 10).Value = "Pera"
$ws.Cells.Item(741, 11).Value = "Packham's Triumph"
$ws.Cells.Item(741, 12).Value = "Segunda"
$ws.Cells.Item(741, 13).Value = 180
$ws.Cells.Item(741, 14).Value = 11000
$ws.Cells.Item(741, 15).Value = 11000
$ws.Cells.Item(741, 16).Value = 11000
$ws.Cells.Item(741, 17).Value = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(741, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(741, 19).Value = 688
$ws.Cells.Item(741, 20).Value = 16
